$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 13.082
$ws.Range("E6").Value = 13.31
$ws.Range("E7").Value = 13.303
$ws.Range("E16").Value = 12.714
$ws.Range("E20").Value = 13.095
